# Update "想去人数" (F column) figures across the workbook's sheets to
# reflect refreshed scrape counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 315
$ws1.Cells.Item(3, 6).Value = 1261
$ws1.Cells.Item(5, 6).Value = 332
$ws1.Cells.Item(6, 6).Value = 3843
$ws1.Cells.Item(9, 6).Value = 2224
$ws1.Cells.Item(10, 6).Value = 335
$ws1.Cells.Item(11, 6).Value = 219
$ws1.Cells.Item(12, 6).Value = 738
$ws1.Cells.Item(13, 6).Value = 160
$ws1.Cells.Item(15, 6).Value = 2129
$ws1.Cells.Item(20, 6).Value = 225
$ws1.Cells.Item(21, 6).Value = 24
$ws1.Cells.Item(22, 6).Value = 268

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(4, 6).Value = 15
$ws2.Cells.Item(9, 6).Value = 94
$ws2.Cells.Item(17, 6).Value = 41

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 6396

# Sheet 4: 全部类型 (All types) - combined view of the other sheets
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 6396
$ws4.Cells.Item(9, 6).Value = 15
$ws4.Cells.Item(10, 6).Value = 315
$ws4.Cells.Item(11, 6).Value = 1261
$ws4.Cells.Item(16, 6).Value = 332
$ws4.Cells.Item(17, 6).Value = 3843
$ws4.Cells.Item(20, 6).Value = 94
$ws4.Cells.Item(24, 6).Value = 2224
$ws4.Cells.Item(25, 6).Value = 335
$ws4.Cells.Item(27, 6).Value = 219
$ws4.Cells.Item(28, 6).Value = 738
$ws4.Cells.Item(29, 6).Value = 160
$ws4.Cells.Item(32, 6).Value = 2129
$ws4.Cells.Item(39, 6).Value = 225
$ws4.Cells.Item(40, 6).Value = 24
$ws4.Cells.Item(42, 6).Value = 41
$ws4.Cells.Item(48, 6).Value = 268
